$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 0.77777777777777779
